# WAZIR_JULY_ATTENDANCE_POWER SYSTEM PROTECTION_EED_2024.xlsx
#
# Edits applied on the "MAY_2024" worksheet:
#   1. D24 (row for roll no. 21ME20) attendance value corrected from 0 to 3
#      -> recalculates the dependent SUM (G24) and percentage (H24) formulas.
#   2. The entire row 25 (roll no. 21ME21 entry) is deleted, shifting the
#      "Teacher" / "Chairman" signature block and the trailing spacer rows
#      up by one row.
#   3. The sheet's print area is adjusted to match the new (one row shorter)
#      layout.
#   4. Selection is left on G9 to match the resulting view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAY_2024")
$ws.Activate()

# 1) Fix the attendance figure for the 21ME20 row (was 0, should be 3).
$ws.Range("D24").Value = 3

# 2) Remove the 21ME21 row entirely (shifts everything below it up by one).
$ws.Rows.Item(25).Delete()

# 3) Shrink the print area by one row to match the new sheet extent.
$ws.PageSetup.PrintArea = '$A$1:$H$26'

# 4) Match the resulting selection/active cell.
$ws.Range("G9").Select()
